$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C38's underlying shared string text (Törmäysdemoa -> extended text) first
# so this edits the existing shared-string slot in place.
$ws.Range("C38").Value = "Törmäysdemoa, rigidbody ja  törmäystarkasteluun liittyvien luokkien tutkintaa"

# Update B38: change from a time-of-day numeric value to a text time range
# (this becomes a brand-new shared string entry)
$ws.Range("B38").Value = "16.00-18.00, 18.45-21-45"

# Match B38's style to the other time-range cells (e.g. B37) which wrap text
# using the same number format as the other "s=3" styled cells.
$ws.Range("B38").NumberFormat = $ws.Range("B37").NumberFormat
$ws.Range("B38").WrapText = $ws.Range("B37").WrapText

# Add new value in G38
$ws.Range("G38").Value = 5

# Row 38 grows taller to fit the wrapped text, matching row 37's height
$ws.Rows.Item(38).RowHeight = 43.2

# Update selection on the sheet
$ws.Range("G42").Select()
